# Adds three new worksheets ("Day 4", "Day 5", "Sheet3") at the end of the
# workbook, populates "Day 4" and "Day 5" with data, applies a date number
# format + wrap-text formatting, and makes "Day 5" the active/selected sheet
# (mirroring the author's upload of new daily-challenge sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Day 4": order_id / seller_id / customer_id / order_date
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$day4 = $wb.Worksheets.Add($null, $lastSheet)
$day4.Name = "Day 4"

$day4.Range("A1").Value = "order_id"
$day4.Range("B1").Value = "seller_id"
$day4.Range("C1").Value = "customer_id"
$day4.Range("D1").Value = "order_date"

$day4.Range("A2").Value = 1
$day4.Range("B2").Value = 3
$day4.Range("C2").Value = 5
$day4.Range("D2").Value = 43678

$day4.Range("A3").Value = 2
$day4.Range("B3").Value = 3
$day4.Range("C3").Value = 6
$day4.Range("D3").Value = 43679

$day4.Range("A4").Value = 3
$day4.Range("B4").Value = 7
$day4.Range("C4").Value = 7
$day4.Range("D4").Value = 43678

$day4.Range("A5").Value = 4
$day4.Range("B5").Value = 7
$day4.Range("C5").Value = 6
$day4.Range("D5").Value = 43679

$day4.Range("A6").Value = 5
$day4.Range("B6").Value = 7
$day4.Range("C6").Value = 1
$day4.Range("D6").Value = 43668

$day4.Range("A7").Value = 6
$day4.Range("B7").Value = 4
$day4.Range("C7").Value = 4
$day4.Range("D7").Value = 43667

$day4.Range("A8").Value = 7
$day4.Range("B8").Value = 4
$day4.Range("C8").Value = 4
$day4.Range("D8").Value = 43667

# Date format on D2, then copy just the format down so every date cell
# shares a single style entry (instead of one new style per cell).
$day4.Range("D2").NumberFormat = "mm-dd-yy"
$day4.Range("D2").Copy() | Out-Null
$day4.Range("D3:D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$day4.Columns.Item(4).ColumnWidth = 10.18

$day4.Application.ActiveWindow.Zoom = 177
$day4.Range("C13").Select() | Out-Null

# ---------------------------------------------------------------------
# "Day 5": review_id / review_text
# ---------------------------------------------------------------------
$day5 = $wb.Worksheets.Add($null, $day4)
$day5.Name = "Day 5"

$day5.Range("A1").Value = "review_id"
$day5.Range("B1").Value = "review_text"

$day5.Range("A2").Value = 1
$day5.Range("B2").Value = "Great hotel with excellent service"

$day5.Range("A3").Value = 2
$day5.Range("B3").Value = "The room was clean and spacious, but the staff was unfriendly."

$day5.Range("A4").Value = 3
$day5.Range("B4").Value = "The hotel was lovely, and the staff were incredibly helpful. Our room had a beautiful view of the city."
$day5.Range("B4").WrapText = $true
$day5.Rows.Item(4).RowHeight = 29

$day5.Columns.Item(2).ColumnWidth = 54.09

$day5.Application.ActiveWindow.Zoom = 210
$day5.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------
# "Sheet3": blank placeholder sheet
# ---------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Add($null, $day5)
$sheet3.Name = "Sheet3"

# "Day 5" ends up the active/selected tab (tabSelected moves off "Day 2").
$day5.Activate() | Out-Null

Write-Output "Added Day 4, Day 5, Sheet3"
